$wb = $excel.ActiveWorkbook

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("D5").Value = "2016-02-15 04:16:49"
$wsZh.Range("G5").Value = "2016-02-15 04:17:46"

$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("D5").Value = "2016-02-15 04:17:03"
$wsDe.Range("G5").Value = "2016-02-15 04:18:11"
